# Refresh the scraped listings in "展览" (Exhibitions), "演出"
# (Performances) and "全部类型" (All types). The stale
# "CM03·配音演员孙路路专场见面会" listing (展览 row 2) and the stale
# "CrossingX意次元｜乐队番 ONLY" listing (演出 row 2) have both
# dropped out of the source feed, so every following row slides up one slot
# and the view-count / price figures for the surviving events are refreshed
# to their newly scraped values. "本地生活" (Local life) is untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "展览" sheet: drop the obsolete first listing (shifts everything
#    else up one row, carrying styles along), then rewrite the surviving
#    rows (B:I) with the freshly scraped values.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Rows.Item(2).EntireRow.Delete()

$expoData = @(
    @(2, '''2024-08-17', '南昌·CM03动漫游戏博览会', '怀玉山大道1315号 南昌绿地国际博览中心', '2024.08.17 09:00-08.18 17:00', 5705, 75, 'https://show.bilibili.com/platform/detail.html?id=89295', '//i2.hdslb.com/bfs/openplatform/202408/YhHLfv5y1722849043508.jpeg'),
    @(3, '''2024-08-18', '九江·如梦令国潮动漫节', '十里大道202号 山水国际大酒店(九江火车站快乐城店)', '2024.08.18 11:00-08.18 17:00', 82, 40, 'https://show.bilibili.com/platform/detail.html?id=90126', '//i1.hdslb.com/bfs/openplatform/202407/bs3xfiQc1721988224155.jpeg'),
    @(4, '''2024-08-24', '于都·希佳微夏日文化交流会', '站前南路23号 赣州于都雅好花园酒店(于都站店)', '2024.08.24 10:00-08.24 16:00', 8, 35, 'https://show.bilibili.com/platform/detail.html?id=90606', '//i1.hdslb.com/bfs/openplatform/202408/SLxwBbc31723445459650.jpeg'),
    @(5, '''2024-08-24', '南昌·第四届龙年动漫展——暑假最后的狂欢', '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆', '2024.08.24 10:00-08.25 18:00', 958, 55, 'https://show.bilibili.com/platform/detail.html?id=87135', '//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'),
    @(6, '''2024-08-24', '赣州·第五人格only', '兴国路恒大帝景西门 江西长庚控股有限公司', '2024.08.24 10:00-08.24 18:00', 157, 55, 'https://show.bilibili.com/platform/detail.html?id=89742', '//i0.hdslb.com/bfs/openplatform/202407/Jxx8Wz6I1721644479535.jpeg'),
    @(7, '''2024-09-15', '南昌·Sunflower Garden动漫游戏展', '怀玉山大道1315号 南昌绿地国际博览中心', '2024.09.15 09:00-09.16 18:00', 2581, 65, 'https://show.bilibili.com/platform/detail.html?id=89659', '//i0.hdslb.com/bfs/openplatform/202407/CQCXbg291721632431682.jpeg'),
    @(8, '''2024-09-15', '南昌·第一届哥布林动漫游戏展——开学季&贺中秋', '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆', '2024.09.15 10:00-09.16 18:00', 86, 55, 'https://show.bilibili.com/platform/detail.html?id=89240', '//i2.hdslb.com/bfs/openplatform/202407/pixnzm5p1720496832036.jpeg'),
    @(9, '''2024-09-15', '赣州·卡尼动漫展', '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心', '2024.09.15 09:30-09.16 17:00', 179, 65, 'https://show.bilibili.com/platform/detail.html?id=90642', '//i1.hdslb.com/bfs/openplatform/202408/VcJiaBPn1723530492504.jpeg'),
    @(10, '''2024-09-15', '鹰潭·MZD动漫游戏嘉年华', '南站路锦都金源酒店18楼 锦都金源酒店', '2024.09.15 10:00-09.15 17:00', 15, 39.9, 'https://show.bilibili.com/platform/detail.html?id=90709', '//i1.hdslb.com/bfs/openplatform/202408/SoFGB10B1723606695453.jpeg'),
    @(11, '''2024-09-17', '南昌·Aud中秋动漫嘉年华', '青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK', '2024.09.17 10:00-09.17 17:00', 94, 29.9, 'https://show.bilibili.com/platform/detail.html?id=90329', '//i0.hdslb.com/bfs/openplatform/202408/pbU7Eftp1722660514298.jpeg'),
    @(12, '''2024-10-01', '九江·星梦次元XACD动漫游戏博览会国庆盛典', '九瑞大道与重庆路交汇处西南角 九江国际会展中心', '2024.10.01 10:00-10.02 17:00', 33, 44.9, 'https://show.bilibili.com/platform/detail.html?id=90732', '//i0.hdslb.com/bfs/openplatform/202408/PTmf3umB1722911829186.jpeg'),
    @(13, '''2024-10-02', '南昌·萌卡动漫展', '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆', '2024.10.02 09:00-10.03 17:00', 2420, 65, 'https://show.bilibili.com/platform/detail.html?id=89738', '//i0.hdslb.com/bfs/openplatform/202407/uqTvacSV1721621530709.jpeg'),
    @(14, '''2024-10-03', '江西·JMG（广电）第二届UP动漫游戏博览会', '怀玉山大道1315号 南昌绿地国际博览中心', '2024.10.03 09:00-10.05 18:00', 458, 19.9, 'https://show.bilibili.com/platform/detail.html?id=90599', '//i2.hdslb.com/bfs/openplatform/202408/2LP6dm961723428231240.jpeg')
)
foreach ($row in $expoData) {
    $r = $row[0]
    $wsExpo.Cells.Item($r, 1).Value2 = ($r - 1)
    $wsExpo.Cells.Item($r, 2).Value2 = $row[1]
    $wsExpo.Cells.Item($r, 3).Value2 = $row[2]
    $wsExpo.Cells.Item($r, 4).Value2 = $row[3]
    $wsExpo.Cells.Item($r, 5).Value2 = $row[4]
    $wsExpo.Cells.Item($r, 6).Value2 = $row[5]
    $wsExpo.Cells.Item($r, 7).Value2 = $row[6]
    $wsExpo.Cells.Item($r, 8).Value2 = $row[7]
    $wsExpo.Cells.Item($r, 9).Value2 = $row[8]
}

# ---------------------------------------------------------------------------
# 2) "演出" sheet: same treatment.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Rows.Item(2).EntireRow.Delete()

$showData = @(
    @(2, '''2024-08-24', '南昌·【8月24日】滑稽互动狂欢大作战《欢乐小丑嘉年华》', '象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院', '2024.08.24 14:30-08.24 20:00', 3, 60, 'https://show.bilibili.com/platform/detail.html?id=90177', '//i1.hdslb.com/bfs/openplatform/202408/ed1EQGH71722479594577.jpeg'),
    @(3, '''2024-09-16', '南昌·《梁祝》65周年大型交响音乐会-风兔子交响乐团', '象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院', '2024.09.16 19:30-09.16 20:50', 1, 140, 'https://show.bilibili.com/platform/detail.html?id=90515', '//i2.hdslb.com/bfs/openplatform/202408/muKn0Ygv1723107475651.jpeg')
)
foreach ($row in $showData) {
    $r = $row[0]
    $wsShow.Cells.Item($r, 1).Value2 = ($r - 1)
    $wsShow.Cells.Item($r, 2).Value2 = $row[1]
    $wsShow.Cells.Item($r, 3).Value2 = $row[2]
    $wsShow.Cells.Item($r, 4).Value2 = $row[3]
    $wsShow.Cells.Item($r, 5).Value2 = $row[4]
    $wsShow.Cells.Item($r, 6).Value2 = $row[5]
    $wsShow.Cells.Item($r, 7).Value2 = $row[6]
    $wsShow.Cells.Item($r, 8).Value2 = $row[7]
    $wsShow.Cells.Item($r, 9).Value2 = $row[8]
}

# ---------------------------------------------------------------------------
# 3) "全部类型" sheet: merges both lists above, so it loses
#    both stale rows (originally rows 2 and 4). Delete bottom-up so row
#    indices don't shift out from under the second delete, then rewrite the
#    survivors.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(4).EntireRow.Delete()
$wsAll.Rows.Item(2).EntireRow.Delete()

$allData = @(
    @(2, '''2024-08-17', '南昌·CM03动漫游戏博览会', '怀玉山大道1315号 南昌绿地国际博览中心', '2024.08.17 09:00-08.18 17:00', 5705, 75, 'https://show.bilibili.com/platform/detail.html?id=89295', '//i2.hdslb.com/bfs/openplatform/202408/YhHLfv5y1722849043508.jpeg'),
    @(3, '''2024-08-18', '九江·如梦令国潮动漫节', '十里大道202号 山水国际大酒店(九江火车站快乐城店)', '2024.08.18 11:00-08.18 17:00', 82, 40, 'https://show.bilibili.com/platform/detail.html?id=90126', '//i1.hdslb.com/bfs/openplatform/202407/bs3xfiQc1721988224155.jpeg'),
    @(4, '''2024-08-24', '于都·希佳微夏日文化交流会', '站前南路23号 赣州于都雅好花园酒店(于都站店)', '2024.08.24 10:00-08.24 16:00', 8, 35, 'https://show.bilibili.com/platform/detail.html?id=90606', '//i1.hdslb.com/bfs/openplatform/202408/SLxwBbc31723445459650.jpeg'),
    @(5, '''2024-08-24', '南昌·【8月24日】滑稽互动狂欢大作战《欢乐小丑嘉年华》', '象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院', '2024.08.24 14:30-08.24 20:00', 3, 60, 'https://show.bilibili.com/platform/detail.html?id=90177', '//i1.hdslb.com/bfs/openplatform/202408/ed1EQGH71722479594577.jpeg'),
    @(6, '''2024-08-24', '南昌·第四届龙年动漫展——暑假最后的狂欢', '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆', '2024.08.24 10:00-08.25 18:00', 958, 55, 'https://show.bilibili.com/platform/detail.html?id=87135', '//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'),
    @(7, '''2024-08-24', '赣州·第五人格only', '兴国路恒大帝景西门 江西长庚控股有限公司', '2024.08.24 10:00-08.24 18:00', 157, 55, 'https://show.bilibili.com/platform/detail.html?id=89742', '//i0.hdslb.com/bfs/openplatform/202407/Jxx8Wz6I1721644479535.jpeg'),
    @(8, '''2024-09-15', '南昌·Sunflower Garden动漫游戏展', '怀玉山大道1315号 南昌绿地国际博览中心', '2024.09.15 09:00-09.16 18:00', 2581, 65, 'https://show.bilibili.com/platform/detail.html?id=89659', '//i0.hdslb.com/bfs/openplatform/202407/CQCXbg291721632431682.jpeg'),
    @(9, '''2024-09-15', '南昌·第一届哥布林动漫游戏展——开学季&贺中秋', '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆', '2024.09.15 10:00-09.16 18:00', 86, 55, 'https://show.bilibili.com/platform/detail.html?id=89240', '//i2.hdslb.com/bfs/openplatform/202407/pixnzm5p1720496832036.jpeg'),
    @(10, '''2024-09-15', '赣州·卡尼动漫展', '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心', '2024.09.15 09:30-09.16 17:00', 179, 65, 'https://show.bilibili.com/platform/detail.html?id=90642', '//i1.hdslb.com/bfs/openplatform/202408/VcJiaBPn1723530492504.jpeg'),
    @(11, '''2024-09-15', '鹰潭·MZD动漫游戏嘉年华', '南站路锦都金源酒店18楼 锦都金源酒店', '2024.09.15 10:00-09.15 17:00', 15, 39.9, 'https://show.bilibili.com/platform/detail.html?id=90709', '//i1.hdslb.com/bfs/openplatform/202408/SoFGB10B1723606695453.jpeg'),
    @(12, '''2024-09-16', '南昌·《梁祝》65周年大型交响音乐会-风兔子交响乐团', '象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院', '2024.09.16 19:30-09.16 20:50', 1, 140, 'https://show.bilibili.com/platform/detail.html?id=90515', '//i2.hdslb.com/bfs/openplatform/202408/muKn0Ygv1723107475651.jpeg'),
    @(13, '''2024-09-17', '南昌·Aud中秋动漫嘉年华', '青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK', '2024.09.17 10:00-09.17 17:00', 94, 29.9, 'https://show.bilibili.com/platform/detail.html?id=90329', '//i0.hdslb.com/bfs/openplatform/202408/pbU7Eftp1722660514298.jpeg'),
    @(14, '''2024-10-01', '九江·星梦次元XACD动漫游戏博览会国庆盛典', '九瑞大道与重庆路交汇处西南角 九江国际会展中心', '2024.10.01 10:00-10.02 17:00', 33, 44.9, 'https://show.bilibili.com/platform/detail.html?id=90732', '//i0.hdslb.com/bfs/openplatform/202408/PTmf3umB1722911829186.jpeg'),
    @(15, '''2024-10-02', '南昌·萌卡动漫展', '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆', '2024.10.02 09:00-10.03 17:00', 2420, 65, 'https://show.bilibili.com/platform/detail.html?id=89738', '//i0.hdslb.com/bfs/openplatform/202407/uqTvacSV1721621530709.jpeg'),
    @(16, '''2024-10-03', '江西·JMG（广电）第二届UP动漫游戏博览会', '怀玉山大道1315号 南昌绿地国际博览中心', '2024.10.03 09:00-10.05 18:00', 458, 19.9, 'https://show.bilibili.com/platform/detail.html?id=90599', '//i2.hdslb.com/bfs/openplatform/202408/2LP6dm961723428231240.jpeg')
)
foreach ($row in $allData) {
    $r = $row[0]
    $wsAll.Cells.Item($r, 1).Value2 = ($r - 1)
    $wsAll.Cells.Item($r, 2).Value2 = $row[1]
    $wsAll.Cells.Item($r, 3).Value2 = $row[2]
    $wsAll.Cells.Item($r, 4).Value2 = $row[3]
    $wsAll.Cells.Item($r, 5).Value2 = $row[4]
    $wsAll.Cells.Item($r, 6).Value2 = $row[5]
    $wsAll.Cells.Item($r, 7).Value2 = $row[6]
    $wsAll.Cells.Item($r, 8).Value2 = $row[7]
    $wsAll.Cells.Item($r, 9).Value2 = $row[8]
}
